$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 sensor values from 3 decimal places to 2 decimal places
$ws.Range("C5").Value = 10.21
$ws.Range("D5").Value = 0.33
$ws.Range("E5").Value = 28.6
$ws.Range("F5").Value = 24.01
$ws.Range("G5").Value = 10.87
$ws.Range("H5").Value = 44.68
$ws.Range("I5").Value = 15.91
$ws.Range("J5").Value = 7.32
$ws.Range("K5").Value = 11.33
$ws.Range("L5").Value = 11.6
$ws.Range("M5").Value = 12.04
$ws.Range("N5").Value = 3.39
$ws.Range("O5").Value = 10.29
$ws.Range("P5").Value = 15.34
$ws.Range("Q5").Value = 8.57
$ws.Range("R5").Value = 0.35
$ws.Range("S5").Value = 0.31
$ws.Range("T5").Value = 152.74
$ws.Range("U5").Value = 29.54
$ws.Range("V5").Value = 9.880000000000001
$ws.Range("W5").Value = 20.06
$ws.Range("X5").Value = 10.3
$ws.Range("AB5").Value = 7.36
$ws.Range("AC5").Value = 8.710000000000001
$ws.Range("AD5").Value = 12.66
$ws.Range("AE5").Value = 0.52
$ws.Range("AF5").Value = 40.42
$ws.Range("AG5").Value = 5.39
$ws.Range("AH5").Value = 11.98

# Remove row 6 entirely (data trimmed to 5 rows)
$ws.Rows.Item(6).Delete()
